# Remove column M ("column from alcohol data") from the measurement sheet.
# Deleting the whole column shifts everything that was in N left into M,
# which is exactly the row-by-row value pattern described by the diff.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("M:M").Delete()

# The author's re-save also re-zoomed every sheet in the workbook to 75%.
foreach ($ws in $wb.Worksheets) {
  $ws.Activate()
  $excel.ActiveWindow.Zoom = 75
}

# Restore sheet 1 as the active/selected sheet, with the selection parked
# on M1 (the cell that used to hold the now-removed value from N1 in the
# old layout, i.e. right after the deleted column).
$ws1.Activate()
$null = $ws1.Range("M1").Select()
